$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# --- Update "Recommandations" sheet rows 2-44 (columns A-G) ---
# Row 2: BRVM - SERVICES PUBLICS
$wsReco.Range("A2").Value = "BRVM - SERVICES PUBLICS"
$wsReco.Range("B2").Value = 0
$wsReco.Range("C2").Value = 6
$wsReco.Range("D2").Value = 2441.75
$wsReco.Range("E2").Value = 98.89
$wsReco.Range("F2").Value = "🟡 Observer"
$wsReco.Range("G2").Value = "➖ Neutre"

# Row 3: UNIWAX CI
$wsReco.Range("A3").Value = "UNIWAX CI"
$wsReco.Range("B3").Value = 0
$wsReco.Range("C3").Value = 3
$wsReco.Range("D3").Value = 1995
$wsReco.Range("E3").Value = 710
$wsReco.Range("F3").Value = "🟡 Observer"
$wsReco.Range("G3").Value = "➖ Neutre"

# Row 4: CFAO MOTORS CI
$wsReco.Range("A4").Value = "CFAO MOTORS CI"
$wsReco.Range("B4").Value = 0
$wsReco.Range("C4").Value = 3
$wsReco.Range("D4").Value = 1990
$wsReco.Range("E4").Value = 675
$wsReco.Range("F4").Value = "🟡 Observer"
$wsReco.Range("G4").Value = "➖ Neutre"

# Row 5: BRVM - AUTRES SECTEURS
$wsReco.Range("A5").Value = "BRVM - AUTRES SECTEURS"
$wsReco.Range("B5").Value = 0
$wsReco.Range("C5").Value = 3
$wsReco.Range("D5").Value = 1900.28
$wsReco.Range("E5").Value = 632.35
$wsReco.Range("F5").Value = "🟡 Observer"
$wsReco.Range("G5").Value = "➖ Neutre"

# Row 6: NEI-CEDA CI
$wsReco.Range("A6").Value = "NEI-CEDA CI"
$wsReco.Range("B6").Value = 0
$wsReco.Range("C6").Value = 3
$wsReco.Range("D6").Value = 1750
$wsReco.Range("E6").Value = 595
$wsReco.Range("F6").Value = "🟡 Observer"
$wsReco.Range("G6").Value = "➖ Neutre"

# Row 7: SETAO CI
$wsReco.Range("A7").Value = "SETAO CI"
$wsReco.Range("B7").Value = 0
$wsReco.Range("C7").Value = 3
$wsReco.Range("D7").Value = 1735
$wsReco.Range("E7").Value = 580
$wsReco.Range("F7").Value = "🟡 Observer"
$wsReco.Range("G7").Value = "➖ Neutre"

# Row 8: AIR LIQUIDE CI
$wsReco.Range("A8").Value = "AIR LIQUIDE CI"
$wsReco.Range("B8").Value = 0
$wsReco.Range("C8").Value = 3
$wsReco.Range("D8").Value = 1640
$wsReco.Range("E8").Value = 550
$wsReco.Range("F8").Value = "🟡 Observer"
$wsReco.Range("G8").Value = "➖ Neutre"

# Row 9: BRVM - DISTRIBUTION
$wsReco.Range("A9").Value = "BRVM - DISTRIBUTION"
$wsReco.Range("B9").Value = 0
$wsReco.Range("C9").Value = 3
$wsReco.Range("D9").Value = 1118.99
$wsReco.Range("E9").Value = 374.01
$wsReco.Range("F9").Value = "🟡 Observer"
$wsReco.Range("G9").Value = "➖ Neutre"

# Row 10: BRVM - TRANSPORT
$wsReco.Range("A10").Value = "BRVM - TRANSPORT"
$wsReco.Range("B10").Value = 0
$wsReco.Range("C10").Value = 3
$wsReco.Range("D10").Value = 1107.6
$wsReco.Range("E10").Value = 373.28
$wsReco.Range("F10").Value = "🟡 Observer"
$wsReco.Range("G10").Value = "➖ Neutre"

# Row 11: SAFCA CI
$wsReco.Range("A11").Value = "SAFCA CI"
$wsReco.Range("B11").Value = 0
$wsReco.Range("C11").Value = 1
$wsReco.Range("D11").Value = 995
$wsReco.Range("E11").Value = 995
$wsReco.Range("F11").Value = "🟡 Observer"
$wsReco.Range("G11").Value = "➖ Neutre"

# Row 12: BRVM - AGRICULTURE
$wsReco.Range("A12").Value = "BRVM - AGRICULTURE"
$wsReco.Range("B12").Value = 0
$wsReco.Range("C12").Value = 3
$wsReco.Range("D12").Value = 974.23
$wsReco.Range("E12").Value = 327.75
$wsReco.Range("F12").Value = "🟡 Observer"
$wsReco.Range("G12").Value = "➖ Neutre"

# Row 13: BRVM - INDUSTRIE
$wsReco.Range("A13").Value = "BRVM - INDUSTRIE"
$wsReco.Range("B13").Value = 0
$wsReco.Range("C13").Value = 3
$wsReco.Range("D13").Value = 796.58
$wsReco.Range("E13").Value = 265.37
$wsReco.Range("F13").Value = "🟡 Observer"
$wsReco.Range("G13").Value = "➖ Neutre"

# Row 14: BRVM - CONSOMMATION DE BASE
$wsReco.Range("A14").Value = "BRVM - CONSOMMATION DE BASE"
$wsReco.Range("B14").Value = 0
$wsReco.Range("C14").Value = 3
$wsReco.Range("D14").Value = 657.09
$wsReco.Range("E14").Value = 219.26
$wsReco.Range("F14").Value = "🟡 Observer"
$wsReco.Range("G14").Value = "➖ Neutre"

# Row 15: BRVM-PRINCIPAL
$wsReco.Range("A15").Value = "BRVM-PRINCIPAL"
$wsReco.Range("B15").Value = 0
$wsReco.Range("C15").Value = 3
$wsReco.Range("D15").Value = 574.47
$wsReco.Range("E15").Value = 192.12
$wsReco.Range("F15").Value = "🟡 Observer"
$wsReco.Range("G15").Value = "➖ Neutre"

# Row 16: BRVM - INDUSTRIELS
$wsReco.Range("A16").Value = "BRVM - INDUSTRIELS"
$wsReco.Range("B16").Value = 0
$wsReco.Range("C16").Value = 3
$wsReco.Range("D16").Value = 420.44
$wsReco.Range("E16").Value = 141.27
$wsReco.Range("F16").Value = "🟡 Observer"
$wsReco.Range("G16").Value = "➖ Neutre"

# Row 17: BRVM-PRESTIGE
$wsReco.Range("A17").Value = "BRVM-PRESTIGE"
$wsReco.Range("B17").Value = 0
$wsReco.Range("C17").Value = 3
$wsReco.Range("D17").Value = 393.92
$wsReco.Range("E17").Value = 131.71
$wsReco.Range("F17").Value = "🟡 Observer"
$wsReco.Range("G17").Value = "➖ Neutre"

# Row 18: BRVM - FINANCES
$wsReco.Range("A18").Value = "BRVM - FINANCES"
$wsReco.Range("B18").Value = 0
$wsReco.Range("C18").Value = 3
$wsReco.Range("D18").Value = 371.63
$wsReco.Range("E18").Value = 124.63
$wsReco.Range("F18").Value = "🟡 Observer"
$wsReco.Range("G18").Value = "➖ Neutre"

# Row 19: BRVM - SERVICES FINANCIERS
$wsReco.Range("A19").Value = "BRVM - SERVICES FINANCIERS"
$wsReco.Range("B19").Value = 0
$wsReco.Range("C19").Value = 3
$wsReco.Range("D19").Value = 365.24
$wsReco.Range("E19").Value = 122.49
$wsReco.Range("F19").Value = "🟡 Observer"
$wsReco.Range("G19").Value = "➖ Neutre"

# Row 20: BRVM - ENERGIE
$wsReco.Range("A20").Value = "BRVM - ENERGIE"
$wsReco.Range("B20").Value = 0
$wsReco.Range("C20").Value = 3
$wsReco.Range("D20").Value = 331.78
$wsReco.Range("E20").Value = 110.45
$wsReco.Range("F20").Value = "🟡 Observer"
$wsReco.Range("G20").Value = "➖ Neutre"

# Row 21: BRVM - CONSOMMATION DISCRETIONNAIRE
$wsReco.Range("A21").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsReco.Range("B21").Value = 0
$wsReco.Range("C21").Value = 3
$wsReco.Range("D21").Value = 322.16
$wsReco.Range("E21").Value = 108.39
$wsReco.Range("F21").Value = "🟡 Observer"
$wsReco.Range("G21").Value = "➖ Neutre"

# Row 22: BRVM - TELECOMMUNICATIONS
$wsReco.Range("A22").Value = "BRVM - TELECOMMUNICATIONS"
$wsReco.Range("B22").Value = 0
$wsReco.Range("C22").Value = 3
$wsReco.Range("D22").Value = 279.84
$wsReco.Range("E22").Value = 92.94
$wsReco.Range("F22").Value = "🟡 Observer"
$wsReco.Range("G22").Value = "➖ Neutre"

# Row 23: BERNABE CI (BNBC)
$wsReco.Range("A23").Value = "BERNABE CI (BNBC)"
$wsReco.Range("B23").Value = 3
$wsReco.Range("C23").Value = 0
$wsReco.Range("D23").Value = 21
$wsReco.Range("E23").Value = 6.72
$wsReco.Range("F23").Value = "🟢 Achat"
$wsReco.Range("G23").Value = "✅ Renforcer"

# Row 24: UNIWAX CI (UNXC)
$wsReco.Range("A24").Value = "UNIWAX CI (UNXC)"
$wsReco.Range("B24").Value = 3
$wsReco.Range("C24").Value = 0
$wsReco.Range("D24").Value = 20.93
$wsReco.Range("E24").Value = 6.77
$wsReco.Range("F24").Value = "🟢 Achat"
$wsReco.Range("G24").Value = "✅ Renforcer"

# Row 25: SUCRIVOIRE (SCRC)
$wsReco.Range("A25").Value = "SUCRIVOIRE (SCRC)"
$wsReco.Range("B25").Value = 2
$wsReco.Range("C25").Value = 0
$wsReco.Range("D25").Value = 14
$wsReco.Range("E25").Value = 7.24
$wsReco.Range("F25").Value = "🟡 Observer"
$wsReco.Range("G25").Value = "➖ Neutre"

# Row 26: AIR LIQUIDE CI (SIVC)
$wsReco.Range("A26").Value = "AIR LIQUIDE CI (SIVC)"
$wsReco.Range("B26").Value = 1
$wsReco.Range("C26").Value = 0
$wsReco.Range("D26").Value = 7.27
$wsReco.Range("E26").Value = 7.27
$wsReco.Range("F26").Value = "🟡 Observer"
$wsReco.Range("G26").Value = "➖ Neutre"

# Row 27: SAFCA CI (SAFC)
$wsReco.Range("A27").Value = "SAFCA CI (SAFC)"
$wsReco.Range("B27").Value = 2
$wsReco.Range("C27").Value = 1
$wsReco.Range("D27").Value = 6.64
$wsReco.Range("E27").Value = 7.04
$wsReco.Range("F27").Value = "🟡 Observer"
$wsReco.Range("G27").Value = "👀 À surveiller"

# Row 28: SETAO CI (STAC)
$wsReco.Range("A28").Value = "SETAO CI (STAC)"
$wsReco.Range("B28").Value = 1
$wsReco.Range("C28").Value = 0
$wsReco.Range("D28").Value = 6.25
$wsReco.Range("E28").Value = 6.25
$wsReco.Range("F28").Value = "🟡 Observer"
$wsReco.Range("G28").Value = "➖ Neutre"

# Row 29: ECOBANK TRANS. INCORP. TG (ETIT)
$wsReco.Range("A29").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$wsReco.Range("B29").Value = 1
$wsReco.Range("C29").Value = 0
$wsReco.Range("D29").Value = 6.25
$wsReco.Range("E29").Value = 6.25
$wsReco.Range("F29").Value = "🟡 Observer"
$wsReco.Range("G29").Value = "➖ Neutre"

# Row 30: AFRICA GLOBAL LOGISTICS CI (SDSC)
$wsReco.Range("A30").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$wsReco.Range("B30").Value = 1
$wsReco.Range("C30").Value = 0
$wsReco.Range("D30").Value = 3.81
$wsReco.Range("E30").Value = 3.81
$wsReco.Range("F30").Value = "🟡 Observer"
$wsReco.Range("G30").Value = "➖ Neutre"

# Row 31: ORAGROUP TOGO (ORGT)
$wsReco.Range("A31").Value = "ORAGROUP TOGO (ORGT)"
$wsReco.Range("B31").Value = 1
$wsReco.Range("C31").Value = 1
$wsReco.Range("D31").Value = 1.89
$wsReco.Range("E31").Value = 7.26
$wsReco.Range("F31").Value = "🟡 Observer"
$wsReco.Range("G31").Value = "👀 À surveiller"

# Row 32: TOTAL
$wsReco.Range("A32").Value = "TOTAL"
$wsReco.Range("B32").Value = 0
$wsReco.Range("C32").Value = 3
$wsReco.Range("D32").Value = 0
$wsReco.Range("E32").Value = 0
$wsReco.Range("F32").Value = "🟡 Observer"
$wsReco.Range("G32").Value = "➖ Neutre"

# Row 33: TOTALENERGIES MARKETING SN (TTLS)
$wsReco.Range("A33").Value = "TOTALENERGIES MARKETING SN (TTLS)"
$wsReco.Range("B33").Value = 0
$wsReco.Range("C33").Value = 1
$wsReco.Range("D33").Value = -0.2
$wsReco.Range("E33").Value = -0.2
$wsReco.Range("F33").Value = "🟡 Observer"
$wsReco.Range("G33").Value = "➖ Neutre"

# Row 34: ECOBANK COTE D''IVOIRE (ECOC)
$wsReco.Range("A34").Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$wsReco.Range("B34").Value = 0
$wsReco.Range("C34").Value = 1
$wsReco.Range("D34").Value = -0.42
$wsReco.Range("E34").Value = -0.42
$wsReco.Range("F34").Value = "🟡 Observer"
$wsReco.Range("G34").Value = "➖ Neutre"

# Row 35: LOTERIE NATIONALE DU BENIN (LNBB)
$wsReco.Range("A35").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$wsReco.Range("B35").Value = 0
$wsReco.Range("C35").Value = 1
$wsReco.Range("D35").Value = -0.55
$wsReco.Range("E35").Value = -0.55
$wsReco.Range("F35").Value = "🟡 Observer"
$wsReco.Range("G35").Value = "➖ Neutre"

# Row 36: NSIA BANQUE COTE D'IVOIRE (NSBC)
$wsReco.Range("A36").Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$wsReco.Range("B36").Value = 0
$wsReco.Range("C36").Value = 1
$wsReco.Range("D36").Value = -0.89
$wsReco.Range("E36").Value = -0.89
$wsReco.Range("F36").Value = "🟡 Observer"
$wsReco.Range("G36").Value = "➖ Neutre"

# Row 37: ORANGE COTE D'IVOIRE (ORAC)
$wsReco.Range("A37").Value = "ORANGE COTE D'IVOIRE (ORAC)"
$wsReco.Range("B37").Value = 0
$wsReco.Range("C37").Value = 1
$wsReco.Range("D37").Value = -1.7
$wsReco.Range("E37").Value = -1.7
$wsReco.Range("F37").Value = "🟡 Observer"
$wsReco.Range("G37").Value = "➖ Neutre"

# Row 38: CIE CI (CIEC)
$wsReco.Range("A38").Value = "CIE CI (CIEC)"
$wsReco.Range("B38").Value = 0
$wsReco.Range("C38").Value = 1
$wsReco.Range("D38").Value = -1.86
$wsReco.Range("E38").Value = -1.86
$wsReco.Range("F38").Value = "🟡 Observer"
$wsReco.Range("G38").Value = "➖ Neutre"

# Row 39: SICABLE CI (CABC)
$wsReco.Range("A39").Value = "SICABLE CI (CABC)"
$wsReco.Range("B39").Value = 0
$wsReco.Range("C39").Value = 1
$wsReco.Range("D39").Value = -2.25
$wsReco.Range("E39").Value = -2.25
$wsReco.Range("F39").Value = "🟡 Observer"
$wsReco.Range("G39").Value = "➖ Neutre"

# Row 40: VIVO ENERGY CI (SHEC)
$wsReco.Range("A40").Value = "VIVO ENERGY CI (SHEC)"
$wsReco.Range("B40").Value = 0
$wsReco.Range("C40").Value = 1
$wsReco.Range("D40").Value = -2.53
$wsReco.Range("E40").Value = -2.53
$wsReco.Range("F40").Value = "🟡 Observer"
$wsReco.Range("G40").Value = "➖ Neutre"

# Row 41: BICI CI (BICC)
$wsReco.Range("A41").Value = "BICI CI (BICC)"
$wsReco.Range("B41").Value = 0
$wsReco.Range("C41").Value = 2
$wsReco.Range("D41").Value = -2.73
$wsReco.Range("E41").Value = -0.92
$wsReco.Range("F41").Value = "🟡 Observer"
$wsReco.Range("G41").Value = "➖ Neutre"

# Row 42: SOGB CI (SOGC)
$wsReco.Range("A42").Value = "SOGB CI (SOGC)"
$wsReco.Range("B42").Value = 0
$wsReco.Range("C42").Value = 1
$wsReco.Range("D42").Value = -3.16
$wsReco.Range("E42").Value = -3.16
$wsReco.Range("F42").Value = "🟡 Observer"
$wsReco.Range("G42").Value = "➖ Neutre"

# Row 43: SERVAIR ABIDJAN CI (ABJC)
$wsReco.Range("A43").Value = "SERVAIR ABIDJAN CI (ABJC)"
$wsReco.Range("B43").Value = 0
$wsReco.Range("C43").Value = 1
$wsReco.Range("D43").Value = -3.68
$wsReco.Range("E43").Value = -3.68
$wsReco.Range("F43").Value = "🟡 Observer"
$wsReco.Range("G43").Value = "➖ Neutre"

# Row 44: SOLIBRA CI (SLBC)
$wsReco.Range("A44").Value = "SOLIBRA CI (SLBC)"
$wsReco.Range("B44").Value = 0
$wsReco.Range("C44").Value = 1
$wsReco.Range("D44").Value = -5.06
$wsReco.Range("E44").Value = -5.06
$wsReco.Range("F44").Value = "🟡 Observer"
$wsReco.Range("G44").Value = "👀 À surveiller"

# --- Update "Top_YTD" sheet rows 2-11 (columns A-B) ---
# Row 2: BRVM - SERVICES PUBLICS
$wsYtd.Range("A2").Value = "BRVM - SERVICES PUBLICS"
$wsYtd.Range("B2").Value = 429759.24

# Row 3: UNIWAX CI
$wsYtd.Range("A3").Value = "UNIWAX CI"
$wsYtd.Range("B3").Value = 44514.8

# Row 4: CFAO MOTORS CI
$wsYtd.Range("A4").Value = "CFAO MOTORS CI"
$wsYtd.Range("B4").Value = 44365.62

# Row 5: BRVM - AUTRES SECTEURS
$wsYtd.Range("A5").Value = "BRVM - AUTRES SECTEURS"
$wsYtd.Range("B5").Value = 39351.84

# Row 6: NEI-CEDA CI
$wsYtd.Range("A6").Value = "NEI-CEDA CI"
$wsYtd.Range("B6").Value = 31765.75

# Row 7: SETAO CI
$wsYtd.Range("A7").Value = "SETAO CI"
$wsYtd.Range("B7").Value = 31091.6

# Row 8: AIR LIQUIDE CI
$wsYtd.Range("A8").Value = "AIR LIQUIDE CI"
$wsYtd.Range("B8").Value = 26941.62

# Row 9: BRVM - DISTRIBUTION
$wsYtd.Range("A9").Value = "BRVM - DISTRIBUTION"
$wsYtd.Range("B9").Value = 10482.12

# Row 10: BRVM - TRANSPORT
$wsYtd.Range("A10").Value = "BRVM - TRANSPORT"
$wsYtd.Range("B10").Value = 10228.79

# Row 11: BRVM - AGRICULTURE
$wsYtd.Range("A11").Value = "BRVM - AGRICULTURE"
$wsYtd.Range("B11").Value = 7562.36

Write-Host "Update complete"